# ---------------------------------------------------------------
# Weekly CompStat update (cs-en-us-006pct): new crime data collected
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number and week-covering date range ---
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# --- Cells that flip from a text placeholder to a real number ---
# (copy number formatting from a sibling numeric cell in the same row first,
#  then overwrite with the new value, so the style/number format matches)
$ws.Range("D20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("C22").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("H22").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 100

# --- Cells that flip from a real number back to a text placeholder ---
# (copy directly from a sibling cell that already holds the same placeholder text)
$ws.Range("G31").Copy($ws.Range("F31"))

# --- Plain numeric value updates (counts and percent changes) ---
# Row 15
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 50

# Row 16
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -76.923076923076
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = -55.555555555555
$ws.Range("L16").Value = -63.636363636363
$ws.Range("M16").Value = -56.756756756756
$ws.Range("N16").Value = -92.344497607655

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = -12.5
$ws.Range("L17").Value = -48.780487804878
$ws.Range("N17").Value = -55.31914893617

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = -36.842105263157
$ws.Range("L18").Value = -50.684931506849
$ws.Range("M18").Value = -36.842105263157
$ws.Range("N18").Value = -79.428571428571

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = -9.859154929577
$ws.Range("I19").Value = 212
$ws.Range("J19").Value = 238
$ws.Range("K19").Value = -10.924369747899
$ws.Range("L19").Value = -30.491803278688
$ws.Range("M19").Value = -8.225108225108
$ws.Range("N19").Value = -60.299625468164

# Row 20
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = -66.666666666666
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -33.333333333333
$ws.Range("N20").Value = -97.979797979798

# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -15
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -19.642857142857
$ws.Range("I21").Value = 292
$ws.Range("J21").Value = 368
$ws.Range("K21").Value = -20.652173913043
$ws.Range("L21").Value = -38.004246284501
$ws.Range("M21").Value = -17.514124293785
$ws.Range("N21").Value = -74.957118353344

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -7.692307692307
$ws.Range("L22").Value = -7.692307692307
$ws.Range("M22").Value = -40

# Row 24
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -29.729729729729
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = -8.461538461538
$ws.Range("I24").Value = 361
$ws.Range("J24").Value = 399
$ws.Range("K24").Value = -9.523809523809
$ws.Range("L24").Value = -13.429256594724
$ws.Range("M24").Value = 5.247813411078

# Row 25
$ws.Range("C25").Value = 18
$ws.Range("E25").Value = -37.931034482758
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 112
$ws.Range("H25").Value = -19.642857142857
$ws.Range("I25").Value = 268
$ws.Range("J25").Value = 334
$ws.Range("K25").Value = -19.760479041916
$ws.Range("L25").Value = -12.131147540983

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 11.538461538461
$ws.Range("I26").Value = 81
$ws.Range("J26").Value = 82
$ws.Range("K26").Value = -1.219512195121
$ws.Range("L26").Value = -17.34693877551
$ws.Range("M26").Value = 68.75

# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 100

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = 17.647058823529
$ws.Range("L28").Value = 5.263157894736

Write-Output "CompStat weekly figures updated (week of 3/24/2025-3/30/2025, Number 13)."
